$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep a Text number format so Excel does not
# auto-coerce the literal numeric-looking strings (e.g. "304.51") into
# real numbers, which would change both the stored type and precision.
$targetCells = @("D2", "E2", "D3", "E3", "E4", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "E10", "D11", "E11", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "E21", "D22", "E22", "E23", "D24", "E24", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "23.233.10"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "1.604.64"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "304.51"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").Value = "0.3760"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").Value = "52.42"
$ws.Range("E8").Value = "  +5.11%  "
$ws.Range("D9").Value = "0.3624"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").Value = "0.08149"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "22.93"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").Value = "6.604"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").Value = "7.379"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "0.00001248"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "1.604.82"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "93.98"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").Value = "0.06923"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "23.226.92"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("E25").Value = "  +3.23%  "
$ws.Range("D26").Value = "3.072"
$ws.Range("E26").Value = "  +9.94%  "
$ws.Range("D27").Value = "21.18"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").Value = "150.29"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "5.281"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").Value = "135.19"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").Value = "2.391"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").Value = "6.735"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").Value = "1.779.03"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").Value = "0.9625"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").Value = "0.07491"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "10.40"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "0.02767"
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("D38").Value = "0.2521"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "6.123"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").Value = "0.08807"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "1.409"
$ws.Range("E41").Value = "  +3.13%  "
$ws.Range("D42").Value = "0.7097"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("D43").Value = "12.45"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("D44").Value = "15.94"
$ws.Range("E44").Value = "  +4.65%  "
$ws.Range("D45").Value = "0.6538"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "134.12"
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("D49").Value = "0.07943"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("D50").Value = "1.206"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("E51").Value = "  -2.74%  "
